$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(22, 13, 23, 10, 15, 10, 7, 10, 8, 11, 10, 12, 13, 12, 11, 11, 8)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("B19").Formula = "=SUM(B2:B18)"
$ws.Range("B20").Formula = "=B19/60"

$wb.Save()
